# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: row 3 (team "R") ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 551
$wsOff.Range("C3").Value = 391
$wsOff.Range("D3").Value = 129
$wsOff.Range("E3").Value = 73

# --- DEF sheet: row 3 (team "R") ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 530
$wsDef.Range("C3").Value = 341
$wsDef.Range("D3").Value = 117
$wsDef.Range("E3").Value = 55
$wsDef.Range("G3").Value = 13
